$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "1.000" or "0.9990" are not
# auto-converted to numbers by Excel, matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns per latest crypto snapshot
$ws.Range("D2").Value = "29.888.61"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.889.41"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "0.7688"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "242.94"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "25.71"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "0.07177"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("D11").Value = "0.08541"
$ws.Range("E11").Value = "  +5.04%  "
$ws.Range("D12").Value = "0.7642"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "1.910.75"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "5.359"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "93.69"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "6.149"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "29.847.67"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "13.78"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "244.73"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "0.000007816"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "2.149.64"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "0.9990"
$ws.Range("D23").Value = "8.026"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "0.1645"
$ws.Range("E25").Value = "  +3.94%  "
$ws.Range("D26").Value = "9.399"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "162.94"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "2.033"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "1.466"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "4.515"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "4.098"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "0.05449"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "1.244"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "0.7430"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "2.783"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "0.4470"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "1.101.63"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").Value = "73.22"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "6.066"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "0.8532"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "7.667"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").Value = "1.867"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "3.003"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "2.046.53"
$ws.Range("E51").Value = "  +0.12%  "
